# Updated symbol list (Price/Volume refresh + a few coin-row reshuffles in
# rows 15-22, where TigerCash/LEO/GateToken/... shifted down one row and
# CoinExToken moved up to row 15). Columns D/E hold numeric-looking text
# (e.g. "298.21", "1.85%") that must stay text (matches the source
# inlineStr cells), so NumberFormat is forced to "@" before each write to
# stop Excel's COM layer from auto-coercing them to numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '298.21'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.85%'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '3.28%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.013'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.56%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07539'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.85%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.605'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '4.55%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9211'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.98%'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.36%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1184'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.08%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1830'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '5.07%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08985'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '3.31%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04037'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-7.00%'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.45%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001277'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.72%'
$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.04093'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '3.85%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005876'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.72%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.340'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.11%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.374'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.09%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3329'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.21%'
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.280'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '3.88%'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1371'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.10%'
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.3220'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '17.40%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.36%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.003896'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2.94%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001251'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-2.33%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02407'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '4.65%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05223'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '3.06%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.006304'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.28%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007830'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.15%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1327'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.12%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007398'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.26%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007097'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-13.58%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.2970'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.57%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006581'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '5.10%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.04%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.04517'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '33.83%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004202'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.02%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.04%'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.04%'
